{"js": "// 1) Update \"Curso (semestre ideal)\" line to add the EF (7) entry.\nconst courseResults = context.document.body.search(\"Curso (semestre ideal): EP (3)\", { matchCase: true });\ncourseResults.load(\"text\");\nawait context.sync();\nif (courseResults.items.length > 0) {\n  courseResults.items[0].insertText(\"Curso (semestre ideal): EF (7), EP (3)\", Word.InsertLocation.replace);\n}\n\n// 2) Fix the accented \"i\" in the professor's first name.\nconst nameResults = context.document.body.search(\"Fabr\u00edcio Maciel Gomes\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Fabricio Maciel Gomes\", Word.InsertLocation.replace);\n}\n\n// 3) Remove the trailing \"Requisitos\" heading paragraph and the bullet\n//    paragraph that follows it (\"LOQ4203 - Sistemas Produtivos I (Requisito fraco)\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nlet reqIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.replace(/[\\r\\u000b]+$/, \"\");\n  if (paragraphs.items[i].style === \"Heading 2\" && text === \"Requisitos\") {\n    reqIndex = i;\n    break;\n  }\n}\n\nif (reqIndex !== -1) {\n  // Delete from the end backwards so earlier indices stay valid.\n  for (let i = paragraphs.items.length - 1; i >= reqIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update \"Curso (semestre ideal)\" line to add the EF (7) entry.\n$find = $d.Content.Find\n$find.Execute(\"Curso (semestre ideal): EP (3)\", $false, $false, $false, $false, $false, $true, 1, $false, \"Curso (semestre ideal): EF (7), EP (3)\", 2)\n\n# 2) Fix the accented \"i\" in the professor's first name.\n$find2 = $d.Content.Find\n$find2.Execute(\"Fabr\u00edcio Maciel Gomes\", $false, $false, $false, $false, $false, $true, 1, $false, \"Fabricio Maciel Gomes\", 2)\n\n# 3) Remove the trailing \"Requisitos\" heading paragraph and everything after it\n#    (the \"LOQ4203 - Sistemas Produtivos I (Requisito fraco)\" bullet paragraph).\n$count = $d.Paragraphs.Count\n$startIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq \"Requisitos\") {\n        $startIdx = $i\n        break\n    }\n}\n\nif ($startIdx -ge 1) {\n    $pStart = $d.Paragraphs.Item($startIdx)\n    $pEnd = $d.Paragraphs.Item($count)\n    $rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)\n    $rng.Delete()\n}\n\nWrite-Output \"done\"\n"}
